# Update NATMI LR-pair TPM metrics (Fn1-Itgav) with recomputed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "G2" = 29.20950566666667
    "H2" = 87.628517
    "I2" = 0.01829497698069002
    "J2" = 0.01840828041918582
    "M2" = 13.89934866666667
    "N2" = 41.69804600000001
    "O2" = 0.04853507553134179
    "P2" = 0.04999273878390351
    "Q2" = 405.9931036419758
    "R2" = 3653.937932777782
    "S2" = 0.0008879480896019494
    "T2" = 0.0009202803544572025
    "G3" = 29.20950566666667
    "H3" = 87.628517
    "I3" = 0.01829497698069002
    "J3" = 0.01840828041918582
    "O3" = 0.245697991654417
    "P3" = 0.253077086664408
    "Q3" = 2055.249509727493
    "R3" = 18497.24558754743
    "S3" = 0.004495039101519328
    "T3" = 0.004658713978989014
    "G4" = 29.20950566666667
    "H4" = 87.628517
    "I4" = 0.01829497698069002
    "J4" = 0.01840828041918582
    "M4" = 82.007665
    "N4" = 246.022995
    "O4" = 0.2863622109480123
    "P4" = 0.2949625822722868
    "Q4" = 2395.403355527602
    "R4" = 21558.63019974842
    "S4" = 0.005238990057433385
    "T4" = 0.005429753927635423
    "G5" = 29.20950566666667
    "H5" = 87.628517
    "I5" = 0.01829497698069002
    "J5" = 0.01840828041918582
    "M5" = 25.0501465
    "N5" = 50.100293
    "O5" = 0.0874724982879541
    "P5" = 0.06006638442832619
    "Q5" = 731.7023961425803
    "R5" = 4390.214376855482
    "S5" = 0.001600307342621567
    "T5" = 0.001105718848323245
    "G6" = 29.20950566666667
    "H6" = 87.628517
    "I6" = 0.01829497698069002
    "J6" = 0.01840828041918582
    "M6" = 95.05788666666668
    "N6" = 285.17366
    "O6" = 0.3319322235782747
    "P6" = 0.3419012078510756
    "Q6" = 2776.593879251358
    "R6" = 24989.34491326222
    "S6" = 0.006072692389513789
    "T6" = 0.006293813309780937
    "I7" = 0.913374480506715
    "J7" = 0.9190311407684336
    "M7" = 13.89934866666667
    "N7" = 41.69804600000001
    "O7" = 0.04853507553134179
    "P7" = 0.04999273878390351
    "Q7" = 20269.15587375134
    "R7" = 182422.4028637621
    "S7" = 0.04433069939979348
    "T7" = 0.04594488375470916
    "I8" = 0.913374480506715
    "J8" = 0.9190311407684336
    "O8" = 0.245697991654417
    "P8" = 0.253077086664408
    "S8" = 0.2244142754888963
    "T8" = 0.2325857236595426
    "I9" = 0.913374480506715
    "J9" = 0.9190311407684336
    "M9" = 82.007665
    "N9" = 246.022995
    "O9" = 0.2863622109480123
    "P9" = 0.2949625822722868
    "Q9" = 119590.2185484218
    "R9" = 1076311.966935796
    "S9" = 0.2615559356613951
    "T9" = 0.2710797984697026
    "I10" = 0.913374480506715
    "J10" = 0.9190311407684336
    "M10" = 25.0501465
    "N10" = 50.100293
    "O10" = 0.0874724982879541
    "P10" = 0.06006638442832619
    "Q10" = 36530.1523291144
    "R10" = 219180.9139746864
    "S10" = 0.07989514768238459
    "T10" = 0.05520287780299989
    "I11" = 0.913374480506715
    "J11" = 0.9190311407684336
    "M11" = 95.05788666666668
    "N11" = 285.17366
    "O11" = 0.3319322235782747
    "P11" = 0.3419012078510756
    "Q11" = 138621.108663657
    "R11" = 1247589.977972913
    "S11" = 0.3031784222742454
    "T11" = 0.3142178570814794
    "G12" = 57.98602933333333
    "H12" = 173.958088
    "I12" = 0.03631876156896331
    "J12" = 0.03654368891224535
    "M12" = 13.89934866666667
    "N12" = 41.69804600000001
    "O12" = 0.04853507553134179
    "P12" = 0.04999273878390351
    "Q12" = 805.9680394995609
    "R12" = 7253.712355496048
    "S12" = 0.001762733835954428
    "T12" = 0.001826919093990113
    "G13" = 57.98602933333333
    "H13" = 173.958088
    "I13" = 0.03631876156896331
    "J13" = 0.03654368891224535
    "O13" = 0.245697991654417
    "P13" = 0.253077086664408
    "Q13" = 4080.033387705649
    "R13" = 36720.30048935083
    "S13" = 0.00892344677686991
    "T13" = 0.00924837032588148
    "G14" = 57.98602933333333
    "H14" = 173.958088
    "I14" = 0.03631876156896331
    "J14" = 0.03654368891224535
    "M14" = 82.007665
    "N14" = 246.022995
    "O14" = 0.2863622109480123
    "P14" = 0.2949625822722868
    "Q14" = 4755.298868248173
    "R14" = 42797.68981423356
    "S14" = 0.01040032086178204
    "T14" = 0.01077902084731102
    "G15" = 57.98602933333333
    "H15" = 173.958088
    "I15" = 0.03631876156896331
    "J15" = 0.03654368891224535
    "M15" = 25.0501465
    "N15" = 50.100293
    "O15" = 0.0874724982879541
    "P15" = 0.06006638442832619
    "Q15" = 1452.558529753297
    "R15" = 8715.351178519782
    "S15" = 0.003176892809161756
    "T15" = 0.00219504726663209
    "G16" = 57.98602933333333
    "H16" = 173.958088
    "I16" = 0.03631876156896331
    "J16" = 0.03654368891224535
    "M16" = 95.05788666666668
    "N16" = 285.17366
    "O16" = 0.3319322235782747
    "P16" = 0.3419012078510756
    "Q16" = 5512.029404618009
    "R16" = 49608.26464156208
    "S16" = 0.01205536728519518
    "T16" = 0.01249433137843064
    "G17" = 29.481085
    "H17" = 58.96217
    "I17" = 0.01846507700595112
    "J17" = 0.01238628926567028
    "M17" = 13.89934866666667
    "N17" = 41.69804600000001
    "O17" = 0.04853507553134179
    "P17" = 0.04999273878390351
    "Q17" = 409.7678794866367
    "R17" = 2458.60727691982
    "S17" = 0.0008962039071758804
    "T17" = 0.0006192245237605222
    "G18" = 29.481085
    "H18" = 58.96217
    "I18" = 0.01846507700595112
    "J18" = 0.01238628926567028
    "O18" = 0.245697991654417
    "P18" = 0.253077086664408
    "Q18" = 2074.358470284892
    "R18" = 12446.15082170935
    "S18" = 0.004536832336106347
    "T18" = 0.003134686001938463
    "G19" = 29.481085
    "H19" = 58.96217
    "I19" = 0.01846507700595112
    "J19" = 0.01238628926567028
    "M19" = 82.007665
    "N19" = 246.022995
    "O19" = 0.2863622109480123
    "P19" = 0.2949625822722868
    "Q19" = 2417.674942516525
    "R19" = 14506.04965509915
    "S19" = 0.005287700276749468
    "T19" = 0.003653491866573612
    "G20" = 29.481085
    "H20" = 58.96217
    "I20" = 0.01846507700595112
    "J20" = 0.01238628926567028
    "M20" = 25.0501465
    "N20" = 50.100293
    "O20" = 0.0874724982879541
    "P20" = 0.06006638442832619
    "Q20" = 738.5054982289525
    "R20" = 2954.02199291581
    "S20" = 0.00161518641679
    "T20" = 0.0007439996126722009
    "G21" = 29.481085
    "H21" = 58.96217
    "I21" = 0.01846507700595112
    "J21" = 0.01238628926567028
    "M21" = 95.05788666666668
    "N21" = 285.17366
    "O21" = 0.3319322235782747
    "P21" = 0.3419012078510756
    "Q21" = 2802.409636740367
    "R21" = 16814.4578204422
    "S21" = 0.006129154069129427
    "T21" = 0.004234887260725481
    "G22" = 21.628479
    "H22" = 64.885437
    "I22" = 0.01354670393768061
    "J22" = 0.01363060063446486
    "M22" = 13.89934866666667
    "N22" = 41.69804600000001
    "O22" = 0.04853507553134179
    "P22" = 0.04999273878390351
    "Q22" = 300.621770750678
    "R22" = 2705.595936756102
    "S22" = 0.0006574902988160537
    "T22" = 0.0006814310569865112
    "G23" = 21.628479
    "H23" = 64.885437
    "I23" = 0.01354670393768061
    "J23" = 0.01363060063446486
    "O23" = 0.245697991654417
    "P23" = 0.253077086664408
    "Q23" = 1521.830645413115
    "R23" = 13696.47580871803
    "S23" = 0.003328397951025109
    "T23" = 0.003449592698056398
    "G24" = 21.628479
    "H24" = 64.885437
    "I24" = 0.01354670393768061
    "J24" = 0.01363060063446486
    "M24" = 82.007665
    "N24" = 246.022995
    "O24" = 0.2863622109480123
    "P24" = 0.2949625822722868
    "Q24" = 1773.701060291535
    "R24" = 15963.30954262381
    "S24" = 0.003879264090652364
    "T24" = 0.004020517161064026
    "G25" = 21.628479
    "H25" = 64.885437
    "I25" = 0.01354670393768061
    "J25" = 0.01363060063446486
    "M25" = 25.0501465
    "N25" = 50.100293
    "O25" = 0.0874724982879541
    "P25" = 0.06006638442832619
    "Q25" = 541.7965675221735
    "R25" = 3250.779405133041
    "S25" = 0.001184964036996188
    "T25" = 0.0008187408976987532
    "G26" = 21.628479
    "H26" = 64.885437
    "I26" = 0.01354670393768061
    "J26" = 0.01363060063446486
    "M26" = 95.05788666666668
    "N26" = 285.17366
    "O26" = 0.3319322235782747
    "P26" = 0.3419012078510756
    "Q26" = 2055.95750555438
    "R26" = 18503.61754998942
    "S26" = 0.004496587560190894
    "T26" = 0.004660318820659174
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
